$d = $word.ActiveDocument
$replacements = @(
    @("64-36=", "94+4="),
    @("53-35=", "80-14="),
    @("64-32=", "95-82="),
    @("53-15=", "55-20="),
    @("72-58=", "96-83="),
    @("43+1=", "80-39="),
    @("54-48=", "87-30="),
    @("43-0=", "33+11="),
    @("8+57=", "85-75="),
    @("46-11=", "66+17="),
    @("13-0=", "79-34="),
    @("11+55=", "74-30="),
    @("95-50=", "70-34="),
    @("44+13=", "11+18="),
    @("83-69=", "4+1="),
    @("84-22=", "42-12="),
    @("9+53=", "49+18="),
    @("66-60=", "36+47="),
    @("25+5=", "92-7="),
    @("20+46=", "69+23="),
    @("70-48=", "93-83="),
    @("22+60=", "13+56="),
    @("38-28=", "1+70="),
    @("46+2=", "94-74="),
    @("81-76=", "81-50="),
    @("15+18=", "36+5="),
    @("7+92=", "56+19="),
    @("51-7=", "68-22="),
    @("58-24=", "13+79="),
    @("12+11=", "60-17="),
    @("11+69=", "40-34="),
    @("73-38=", "43-36="),
    @("77+3=", "53+37="),
    @("63-20=", "81-13="),
    @("21-3=", "26+9="),
    @("4+13=", "98-38="),
    @("24+31=", "19+62="),
    @("7+79=", "67-7="),
    @("7+48=", "71-17="),
    @("6+56=", "15+82="),
    @("75+23=", "96-83="),
    @("30+50=", "59-10="),
    @("52+36=", "88+11="),
    @("84+6=", "52-14="),
    @("15+43=", "9+6="),
    @("56-42=", "80-25="),
    @("42+6=", "94-45="),
    @("78-2=", "69+11="),
    @("40-19=", "93-14="),
    @("24-23=", "20+34="),
    @("78+21=", "74-60="),
    @("86-46=", "64-57="),
    @("66-27=", "30-29="),
    @("92-73=", "33+44="),
    @("9-7=", "1+78="),
    @("17+33=", "18+4="),
    @("92-71=", "29+46="),
    @("54-29=", "67-46="),
    @("95-9=", "5+27="),
    @("5+37=", "54+25="),
    @("69-54=", "23+76="),
    @("17+3=", "97-31="),
    @("59-6=", "76-26="),
    @("1+87=", "84-14="),
    @("24+45=", "32-0="),
    @("59+26=", "22+47="),
    @("4+84=", "68-14="),
    @("78-30=", "20+29="),
    @("96-55=", "11+34="),
    @("99-55=", "53+32="),
    @("60-11=", "21+66="),
    @("20-1=", "57+26="),
    @("63+24=", "40+25="),
    @("69+0=", "14+53="),
    @("84-84=", "54-27="),
    @("29+1=", "5-0="),
    @("20+61=", "93-0="),
    @("99-80=", "71-63="),
    @("16+83=", "60-37="),
    @("13+9=", "38+22="),
    @("4-0=", "4+22="),
    @("53-14=", "49+14="),
    @("52-26=", "87-75="),
    @("82-81=", "13+29="),
    @("48+44=", "68-0="),
    @("84-12=", "20+29="),
    @("76-37=", "21+22="),
    @("44+17=", "18+61="),
    @("4+66=", "59-23="),
    @("64+0=", "11+86="),
    @("85-44=", "65-43="),
    @("69-52=", "51-17="),
    @("24-1=", "86-14="),
    @("8+36=", "13+34="),
    @("12+74=", "90-20="),
    @("26-25=", "69-7="),
    @("1+74=", "61-49="),
    @("4+83=", "95-48="),
    @("90-39=", "28-19="),
    @("39+24=", "30+26="),
)

foreach ($pair in $replacements) {
    $old = $pair[0]
    $new = $pair[1]
    $r = $d.Content
    $found = $r.Find.Execute($old, $true, $true, $false, $false, $false, $true, 1, $false, $new, 2)
    if (-not $found) {
        Write-Host "NOT FOUND: $old"
    }
}
Write-Host "Done"
